$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 140, shifting existing rows 140-238 down to 141-239
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with its new data
$ws.Cells.Item(140, 1).Value = 10
$ws.Cells.Item(140, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(140, 3).Value = "La Araucanía"
$ws.Cells.Item(140, 4).Value = 45086
$ws.Cells.Item(140, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(140, 5).Value = 9
$ws.Cells.Item(140, 6).Value = "Fruta"
$ws.Cells.Item(140, 7).Value = 100104
$ws.Cells.Item(140, 8).Value = "Frutos de pepita"
$ws.Cells.Item(140, 9).Value = 100104001
$ws.Cells.Item(140, 10).Value = "Granada"
$ws.Cells.Item(140, 11).Value = "Wonderfull"
$ws.Cells.Item(140, 12).Value = "Primera"
$ws.Cells.Item(140, 13).Value = 75
$ws.Cells.Item(140, 14).Value = 13000
$ws.Cells.Item(140, 15).Value = 13000
$ws.Cells.Item(140, 16).Value = 13000
$ws.Cells.Item(140, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(140, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(140, 19).Value = 1300
$ws.Cells.Item(140, 20).Value = 10
